$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster: ECs) -- text value unchanged, numeric stats recomputed with new TPM data
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01237333333333333
$ws.Range("H2").Value = 0.03712
$ws.Range("M2").Value = 4.850457666666666
$ws.Range("N2").Value = 14.551373
$ws.Range("O2").Value = 0.1317659673285644
$ws.Range("P2").Value = 0.1354819285975933
$ws.Range("Q2").Value = 0.06001632952888888
$ws.Range("R2").Value = 0.54014696576
$ws.Range("S2").Value = 0.1317659673285644
$ws.Range("T2").Value = 0.1354819285975933

# Row 3 (Target cluster: FAPs)
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01237333333333333
$ws.Range("H3").Value = 0.03712
$ws.Range("O3").Value = 0.208027493838598
$ws.Range("P3").Value = 0.2138941233307932
$ws.Range("Q3").Value = 0.09475167886222223
$ws.Range("R3").Value = 0.85276510976
$ws.Range("S3").Value = 0.208027493838598
$ws.Range("T3").Value = 0.2138941233307932

# Row 4 (Target cluster: Inflammatory-Mac)
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01237333333333333
$ws.Range("H4").Value = 0.03712
$ws.Range("M4").Value = 10.683974
$ws.Range("N4").Value = 32.051922
$ws.Range("O4").Value = 0.290237389081408
$ws.Range("P4").Value = 0.2984224380626923
$ws.Range("Q4").Value = 0.1321963716266667
$ws.Range("R4").Value = 1.18976734464
$ws.Range("S4").Value = 0.290237389081408
$ws.Range("T4").Value = 0.2984224380626923

# Row 5 (Target cluster: MuSCs)
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.01237333333333333
$ws.Range("H5").Value = 0.03712
$ws.Range("M5").Value = 3.028939
$ws.Range("N5").Value = 6.057878000000001
$ws.Range("O5").Value = 0.08228317918471638
$ws.Range("P5").Value = 0.05640244358033648
$ws.Range("Q5").Value = 0.03747807189333334
$ws.Range("R5").Value = 0.22486843136
$ws.Range("S5").Value = 0.08228317918471638
$ws.Range("T5").Value = 0.05640244358033648

# Row 6 (Target cluster: Resolving-Mac)
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.01237333333333333
$ws.Range("H6").Value = 0.03712
$ws.Range("M6").Value = 10.59005333333333
$ws.Range("N6").Value = 31.77016
$ws.Range("O6").Value = 0.2876859705667131
$ws.Range("P6").Value = 0.2957990664285848
$ws.Range("Q6").Value = 0.1310342599111111
$ws.Range("R6").Value = 1.1793083392
$ws.Range("S6").Value = 0.2876859705667131
$ws.Range("T6").Value = 0.2957990664285848
